# Apply updates to "Unidades Pedido" (col L) and "Diferencia Stock" (col M)
# for a set of product rows, and refresh the corresponding summary metrics
# (Total_Unidades / Total_Ajuste_Stock) so that all sections reconcile.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana_14")

# Row => new L (Unidades Pedido) value, new M (Diferencia Stock) value
$updates = @{
    2  = @(1, 0)
    3  = @(6, 0)
    4  = @(5, 0)
    6  = @(1, 0)
    7  = @(4, 0)
    16 = @(1, 0)
    17 = @(5, 0)
    18 = @(43, 0)
    19 = @(3, 0)
    21 = @(8, 0)
    22 = @(4, 0)
    24 = @(4, 0)
    27 = @(10, 0)
    34 = @(5, 0)
    35 = @(5, 0)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 12).Value = $vals[0]   # Column L = 12
    $ws.Cells.Item($row, 13).Value = $vals[1]   # Column M = 13
}

# Recalculate the summary metrics to reflect the new totals
$totalUnidades = 0
$totalAjusteStock = 0
for ($r = 2; $r -le 37; $r++) {
    $totalUnidades += $ws.Cells.Item($r, 12).Value2
    $totalAjusteStock += $ws.Cells.Item($r, 13).Value2
}

$ws.Cells.Item(40, 3).Value = $totalUnidades       # Total_Unidades
$ws.Cells.Item(51, 3).Value = $totalAjusteStock    # Total_Ajuste_Stock
